$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, $ref, $value)
    $ws.Range($ref).Value = $value
}

function Clear-CellValue {
    param($ws, $ref)
    $ws.Range($ref).ClearContents()
}

$ws = $wb.Worksheets.Item("ALC")

Set-CellValue $ws "H17" 2621.4
Set-CellValue $ws "J17" 2651.75
Set-CellValue $ws "L17" 7955.25
Set-CellValue $ws "N17" -8291.25
Set-CellValue $ws "H43" 2808.1667
Set-CellValue $ws "I43" 2274.5
Set-CellValue $ws "J43" 3075
Set-CellValue $ws "K43" 2274.5
Set-CellValue $ws "L43" 3075
Set-CellValue $ws "M43" -2205.5
Set-CellValue $ws "N43" -3213
Set-CellValue $ws "H107" 1414.7826
Set-CellValue $ws "I107" 985.1667
Set-CellValue $ws "K107" 985.1667
Set-CellValue $ws "M107" 934.8333
Set-CellValue $ws "H137" 7737.5557
Set-CellValue $ws "I137" 991
Set-CellValue $ws "K137" 2973
Set-CellValue $ws "M137" -423
Set-CellValue $ws "H138" 5867.959
Set-CellValue $ws "J138" 6442
Set-CellValue $ws "L138" 19326
Set-CellValue $ws "N138" -29606
Set-CellValue $ws "H141" 7466
Set-CellValue $ws "I141" 7466
Set-CellValue $ws "J141" 0
Set-CellValue $ws "K141" 22398
Set-CellValue $ws "L141" 0
Set-CellValue $ws "M141" -17218
Clear-CellValue $ws "N141"

$ws = $wb.Worksheets.Item("ARM")

Set-CellValue $ws "H45" 3077.2942
Set-CellValue $ws "I45" 1016.625
Set-CellValue $ws "J45" 4909
Set-CellValue $ws "K45" 1016.625
Set-CellValue $ws "L45" 4909
Set-CellValue $ws "M45" -639.625
Set-CellValue $ws "N45" -5663
Set-CellValue $ws "H61" 2675.2727
Set-CellValue $ws "I61" 2564.6191
Set-CellValue $ws "K61" 2564.6191
Set-CellValue $ws "M61" -2352.6191
Set-CellValue $ws "H74" 104358.27
Set-CellValue $ws "I74" 114244.1
Set-CellValue $ws "K74" 114244.1
Set-CellValue $ws "M74" -113370.1
Set-CellValue $ws "H77" 104358.27
Set-CellValue $ws "I77" 114244.1
Set-CellValue $ws "K77" 571220.5
Set-CellValue $ws "M77" -566852.5
Set-CellValue $ws "H122" 1179.2084
Set-CellValue $ws "I122" 1018.2727
Set-CellValue $ws "K122" 3054.8181
Set-CellValue $ws "M122" -604.8181
Set-CellValue $ws "H132" 26167.83
Set-CellValue $ws "I132" 29825.543
Set-CellValue $ws "K132" 89476.629
Set-CellValue $ws "M132" -86946.629
Set-CellValue $ws "H136" 2675.2727
Set-CellValue $ws "I136" 2564.6191
Set-CellValue $ws "K136" 7693.8573
Set-CellValue $ws "M136" -5143.8573

$ws = $wb.Worksheets.Item("BSM")

Set-CellValue $ws "H86" 1966.3572
Set-CellValue $ws "I86" 1566.375
Set-CellValue $ws "K86" 1566.375
Set-CellValue $ws "M86" -443.375
Set-CellValue $ws "H89" 1966.3572
Set-CellValue $ws "I89" 1566.375
Set-CellValue $ws "K89" 7831.875
Set-CellValue $ws "M89" -2215.875
Set-CellValue $ws "H99" 52055.76
Set-CellValue $ws "I99" 80716.16
Set-CellValue $ws "K99" 80716.16
Set-CellValue $ws "M99" -79218.16

$ws = $wb.Worksheets.Item("CRP")

Set-CellValue $ws "H31" 5180.1304
Set-CellValue $ws "I31" 2861.1
Set-CellValue $ws "K31" 2861.1
Set-CellValue $ws "M31" -2566.1
Set-CellValue $ws "H34" 5180.1304
Set-CellValue $ws "I34" 2861.1
Set-CellValue $ws "K34" 2861.1
Set-CellValue $ws "M34" -2659.1
Set-CellValue $ws "H58" 70258.47
Set-CellValue $ws "I58" 74955.5
Set-CellValue $ws "J58" 4500
Set-CellValue $ws "K58" 74955.5
Set-CellValue $ws "L58" 4500
Set-CellValue $ws "M58" -74752.5
Set-CellValue $ws "N58" -4906
Set-CellValue $ws "H107" 507.35715
Set-CellValue $ws "I107" 442.75
Set-CellValue $ws "J107" 895
Set-CellValue $ws "K107" 442.75
Set-CellValue $ws "L107" 895
Set-CellValue $ws "M107" 1477.25
Set-CellValue $ws "N107" -4735
Set-CellValue $ws "H122" 1717.579
Set-CellValue $ws "I122" 1652.3334
Set-CellValue $ws "J122" 1962.25
Set-CellValue $ws "K122" 4957.0002
Set-CellValue $ws "L122" 5886.75
Set-CellValue $ws "M122" -2507.0002
Set-CellValue $ws "N122" -10786.75
Set-CellValue $ws "H123" 70000
Set-CellValue $ws "J123" 70000
Set-CellValue $ws "L123" 70000
Set-CellValue $ws "N123" -79800
Set-CellValue $ws "H134" 47452.363
Set-CellValue $ws "I134" 51197.6
Set-CellValue $ws "K134" 153592.8
Set-CellValue $ws "M134" -151057.8
Set-CellValue $ws "H136" 70258.47
Set-CellValue $ws "I136" 74955.5
Set-CellValue $ws "J136" 4500
Set-CellValue $ws "K136" 224866.5
Set-CellValue $ws "L136" 13500
Set-CellValue $ws "M136" -222316.5
Set-CellValue $ws "N136" -18600

$ws = $wb.Worksheets.Item("CUL")

Set-CellValue $ws "H94" 14013.889
Set-CellValue $ws "J94" 14585.857
Set-CellValue $ws "L94" 43757.571
Set-CellValue $ws "N94" -45109.571
Set-CellValue $ws "H96" 20005.75
Set-CellValue $ws "J96" 23332.666
Set-CellValue $ws "L96" 69997.99800000001
Set-CellValue $ws "N96" -74115.99800000001
Set-CellValue $ws "H98" 3065.6667
Set-CellValue $ws "I98" 0
Set-CellValue $ws "J98" 3065.6667
Set-CellValue $ws "K98" 0
Set-CellValue $ws "L98" 9197.000100000001
Clear-CellValue $ws "M98"
Set-CellValue $ws "N98" -12193.0001
Set-CellValue $ws "H104" 5389.8
Set-CellValue $ws "I104" 15000
Set-CellValue $ws "J104" 2987.25
Set-CellValue $ws "K104" 45000
Set-CellValue $ws "L104" 8961.75
Set-CellValue $ws "M104" -42379
Set-CellValue $ws "N104" -14203.75
Set-CellValue $ws "H105" 8250
Set-CellValue $ws "J105" 8250
Set-CellValue $ws "L105" 24750
Set-CellValue $ws "N105" -29992
Set-CellValue $ws "H121" 866.087
Set-CellValue $ws "I121" 527.8
Set-CellValue $ws "J121" 960.05554
Set-CellValue $ws "K121" 1583.4
Set-CellValue $ws "L121" 2880.16662
Set-CellValue $ws "M121" -273.3999999999999
Set-CellValue $ws "N121" -5500.16662
Set-CellValue $ws "H138" 1785.5714
Set-CellValue $ws "I138" 1250
Set-CellValue $ws "K138" 3750
Set-CellValue $ws "M138" 1390

$ws = $wb.Worksheets.Item("GSM")

Set-CellValue $ws "H102" 2749.739
Set-CellValue $ws "I102" 1962.25
Set-CellValue $ws "K102" 1962.25
Set-CellValue $ws "M102" -340.25
Set-CellValue $ws "H107" 63849.625
Set-CellValue $ws "I107" 77968.766
Set-CellValue $ws "J107" 2666.6667
Set-CellValue $ws "K107" 77968.766
Set-CellValue $ws "L107" 2666.6667
Set-CellValue $ws "M107" -76048.766
Set-CellValue $ws "N107" -6506.6667
Set-CellValue $ws "H113" 3342.5715
Set-CellValue $ws "I113" 2500
Set-CellValue $ws "K113" 2500
Set-CellValue $ws "M113" -330
Set-CellValue $ws "H126" 5876.25
Set-CellValue $ws "I126" 5042.909
Set-CellValue $ws "J126" 7709.6
Set-CellValue $ws "K126" 15128.727
Set-CellValue $ws "L126" 23128.8
Set-CellValue $ws "M126" -12658.727
Set-CellValue $ws "N126" -28068.8
Set-CellValue $ws "H132" 33640.656
Set-CellValue $ws "I132" 37700.035
Set-CellValue $ws "J132" 5225
Set-CellValue $ws "K132" 113100.105
Set-CellValue $ws "L132" 15675
Set-CellValue $ws "M132" -110570.105
Set-CellValue $ws "N132" -20735

$ws = $wb.Worksheets.Item("LTW")

Set-CellValue $ws "H7" 9670.736999999999
Set-CellValue $ws "I7" 12937.333
Set-CellValue $ws "J7" 4070.8572
Set-CellValue $ws "K7" 12937.333
Set-CellValue $ws "L7" 4070.8572
Set-CellValue $ws "M7" -12825.333
Set-CellValue $ws "N7" -4294.8572
Set-CellValue $ws "H40" 3393
Set-CellValue $ws "I40" 2110.2
Set-CellValue $ws "J40" 6600
Set-CellValue $ws "K40" 2110.2
Set-CellValue $ws "L40" 6600
Set-CellValue $ws "M40" -1974.2
Set-CellValue $ws "N40" -6872
Set-CellValue $ws "H126" 9670.736999999999
Set-CellValue $ws "I126" 12937.333
Set-CellValue $ws "J126" 4070.8572
Set-CellValue $ws "K126" 38811.999
Set-CellValue $ws "L126" 12212.5716
Set-CellValue $ws "M126" -36341.999
Set-CellValue $ws "N126" -17152.5716

$ws = $wb.Worksheets.Item("WVR")

Set-CellValue $ws "H107" 656.3
Set-CellValue $ws "J107" 220
Set-CellValue $ws "L107" 660
Set-CellValue $ws "N107" -4500
Set-CellValue $ws "H122" 868.7
Set-CellValue $ws "I122" 868.7
Set-CellValue $ws "J122" 0
Set-CellValue $ws "K122" 2606.1
Set-CellValue $ws "L122" 0
Set-CellValue $ws "M122" -156.1000000000004
Clear-CellValue $ws "N122"
Set-CellValue $ws "H126" 105518.52
Set-CellValue $ws "I126" 116667.82
Set-CellValue $ws "K126" 350003.46
Set-CellValue $ws "M126" -347533.46
Set-CellValue $ws "H132" 52454.81
Set-CellValue $ws "I132" 59064.84
Set-CellValue $ws "J132" 5240.2856
Set-CellValue $ws "K132" 177194.52
Set-CellValue $ws "L132" 15720.8568
Set-CellValue $ws "M132" -174664.52
Set-CellValue $ws "N132" -20780.8568
Set-CellValue $ws "H136" 3980.725
Set-CellValue $ws "I136" 4089.7742
Set-CellValue $ws "J136" 3605.111
Set-CellValue $ws "K136" 12269.3226
Set-CellValue $ws "L136" 10815.333
Set-CellValue $ws "M136" -9719.3226
Set-CellValue $ws "N136" -15915.333
